$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 02.02.2022 09:31"

# Row 8 (Benzina Albert Modřice): delta price becomes a real number, and the
# scraped timestamp becomes a real Excel date/time value formatted like the
# other rows (style already used by E2:E10).
$ws.Range("D8").Value = 0.1
$ws.Range("E8").Value = 44594.38480324074
$ws.Range("E8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
